$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table values after adding ifoCAST full series evaluation.
# Row => Q0..Q1 (A col unchanged), B=ME, C=MAE, D=MSE, E=RMSE, F=SE, G=N

$ws.Range("B2").Value = -0.1349881739519865
$ws.Range("C2").Value = 1.932678250060537
$ws.Range("D2").Value = 16.79600690525843
$ws.Range("E2").Value = 4.09829316975475
$ws.Range("F2").Value = 4.19246076256923
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = 0.1207090144436769
$ws.Range("C3").Value = 1.786104848928261
$ws.Range("D3").Value = 11.46089950717241
$ws.Range("E3").Value = 3.385395029708116
$ws.Range("F3").Value = 3.466791785362177
$ws.Range("G3").Value = 21

$ws.Range("B4").Value = -0.5555707582142663
$ws.Range("C4").Value = 1.050109602367337
$ws.Range("D4").Value = 4.446451018744308
$ws.Range("E4").Value = 2.108660953957347
$ws.Range("F4").Value = 2.087000403843689
$ws.Range("G4").Value = 20

$ws.Range("B5").Value = 0.07663430359571852
$ws.Range("C5").Value = 0.7265601984156437
$ws.Range("D5").Value = 1.828123414480074
$ws.Range("E5").Value = 1.352081141973393
$ws.Range("F5").Value = 1.386898245123047
$ws.Range("G5").Value = 19

$ws.Range("B6").Value = 0.04069694792458084
$ws.Range("C6").Value = 0.747007170374518
$ws.Range("D6").Value = 1.696764412358688
$ws.Range("E6").Value = 1.302599098862996
$ws.Range("F6").Value = 1.339709078915569
$ws.Range("G6").Value = 18

$ws.Range("B7").Value = -0.003405844617402055
$ws.Range("C7").Value = 0.6048043842844364
$ws.Range("D7").Value = 0.6518033129846593
$ws.Range("E7").Value = 0.8073433674618621
$ws.Range("F7").Value = 0.8321830900003587
$ws.Range("G7").Value = 17

$ws.Range("B8").Value = 0.08389527245345252
$ws.Range("C8").Value = 0.588137197430979
$ws.Range("D8").Value = 0.6090342934726032
$ws.Range("E8").Value = 0.7804064924592844
$ws.Range("F8").Value = 0.8013294797905551
$ws.Range("G8").Value = 16

$ws.Range("B9").Value = 0.225243323359858
$ws.Range("C9").Value = 0.4988262046674656
$ws.Range("D9").Value = 0.4060786597870128
$ws.Range("E9").Value = 0.6372430147024075
$ws.Range("F9").Value = 0.6170298427624463
$ws.Range("G9").Value = 15

$ws.Range("B10").Value = 0.2001121816970861
$ws.Range("C10").Value = 0.4652087079175589
$ws.Range("D10").Value = 0.3905605582328028
$ws.Range("E10").Value = 0.6249484444598633
$ws.Range("F10").Value = 0.614392722160502
$ws.Range("G10").Value = 14

$ws.Range("B11").Value = 0.2336685148374089
$ws.Range("C11").Value = 0.3914692457684976
$ws.Range("D11").Value = 0.1986247317491
$ws.Range("E11").Value = 0.4456733464647622
$ws.Range("F11").Value = 0.395000932487943
